$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (municipio-nombre) moves from a "measure" role to a "dimension" role
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# Column D (energia-renovable) moves from a "dimension" role to a "measure" role
$ws.Range("D2").Value = "iaest-measure:energia-renovable"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"

# The mapping file reference for energia-renovable (row 5) is no longer needed
$ws.Range("D5").Clear()
